$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update manager credentials (rows 2-6, columns A/B)
$ws.Range("A2").Value = "mngr177009"
$ws.Range("B2").Value = "hubapEv"

$ws.Range("A3").Value = "mngr167936"
$ws.Range("B3").Value = "EnutAje"

$ws.Range("A4").Value = "mngr177009"
$ws.Range("B4").Value = "hubapEv"

$ws.Range("A5").Value = "mngr167936"
$ws.Range("B5").Value = "EnutAje"

$ws.Range("A6").Value = "mngr177009"
$ws.Range("B6").Value = "hubapEv"

# Update the active selection to A6:B6 with A6 as the active cell
$ws.Range("A6:B6").Select()
